# Correção das notas do fórum para matc65 em 2021.2
# Zera todos os valores de visualização/nota (colunas B até J, linhas 2 a 50)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B2:J50")
$range.Value = 0
